# Commit: "LGDs added to preprocess"
# Adds two new computed columns (LGDr, LGDc) to the "Input" sheet:
#   N = K / J   (LGDr)
#   O = M / L   (LGDc)
# and restores the various sheet-view bits (zoom, active sheet/tab,
# selection) that Excel/Calc re-wrote when the workbook was re-saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Input sheet: new LGDr / LGDc columns
# ---------------------------------------------------------------------
$input = $wb.Worksheets.Item("Input")

$input.Range("N1").Value2 = "LGDr"
$input.Range("O1").Value2 = "LGDc"

for ($r = 2; $r -le 105; $r++) {
    $input.Range("N$r").Formula = "=K$r/J$r"
    $input.Range("O$r").Formula = "=M$r/L$r"
}

# ---------------------------------------------------------------------
# 2) Sheet view bits: zoom 90 -> 120 on every sheet, updated selections,
#    and the active sheet moving from "I comp" back to "Input".
# ---------------------------------------------------------------------

$hpir = $wb.Worksheets.Item("HPIr comp")
$hpir.Activate() | Out-Null
$hpir.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 120

$hpic = $wb.Worksheets.Item("HPIc comp")
$hpic.Activate() | Out-Null
$hpic.Range("B5").Select() | Out-Null
$excel.ActiveWindow.Zoom = 120

$rhos = $wb.Worksheets.Item("rhos computation")
$rhos.Activate() | Out-Null
$rhos.Range("B15").Select() | Out-Null
$excel.ActiveWindow.Zoom = 120

$icomp = $wb.Worksheets.Item("I comp")
$icomp.Activate() | Out-Null
$icomp.Range("H1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 120

$esc = $wb.Worksheets.Item("Escsount")
$esc.Activate() | Out-Null
$esc.Range("B2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 120

$output = $wb.Worksheets.Item("Output")
$output.Activate() | Out-Null
$output.Range("J2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 120

# Input becomes the active / selected tab again (activeTab 4 -> 0).
$input.Activate() | Out-Null
$input.Range("O2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 120

Write-Output "edit complete"
